$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-11 with newly scraped opportunities.
# NOTE: Column A (OPPORTUNITY ID) values are numeric-looking but must stay
# text, matching the source data's formatting, so they are entered with a
# leading apostrophe (standard Excel "force text" prefix), e.g. '''1327434'
# (a doubled '' inside a single-quoted PowerShell string is a literal ').
# Row 2
$ws.Range("A2").Value = '''1327434'
$ws.Range("B2").Value = 'https://aiesec.org/opportunity/global-talent/1327434'
$ws.Range("C2").Value = 'Internship Software Development with Rust, C++ and Linux (m/f/d)'
$ws.Range("D2").Value = '72 Tübingen, Germany'
$ws.Range("E2").Value = 'No'
$ws.Range("F2").Value = '2 applicants'
$ws.Range("G2").Value = '3 - 6 Months'
$ws.Range("H2").Value = 'Intra2net AG'

# Row 3
$ws.Range("A3").Value = '''1327409'
$ws.Range("B3").Value = 'https://aiesec.org/opportunity/global-talent/1327409'
$ws.Range("C3").Value = 'Web Developer'
$ws.Range("D3").Value = 'Santarém, Portugal'
$ws.Range("E3").Value = 'No'
$ws.Range("F3").Value = '3 applicants'
$ws.Range("G3").Value = '9 - 12 Weeks'
$ws.Range("H3").Value = 'Horas Inversas'

# Row 4
$ws.Range("A4").Value = '''1327406'
$ws.Range("B4").Value = 'https://aiesec.org/opportunity/global-talent/1327406'
$ws.Range("C4").Value = 'Digital Marketing'
$ws.Range("D4").Value = 'Santarém, Portugal'
$ws.Range("E4").Value = 'No'
$ws.Range("F4").Value = '4 applicants'
$ws.Range("G4").Value = '9 - 12 Weeks'
$ws.Range("H4").Value = 'Horas Inversas'

# Row 5
$ws.Range("A5").Value = '''1327381'
$ws.Range("B5").Value = 'https://aiesec.org/opportunity/global-talent/1327381'
$ws.Range("C5").Value = 'Product Management Intern'
$ws.Range("D5").Value = 'Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye'
$ws.Range("E5").Value = 'No'
$ws.Range("F5").Value = '6 applicants'
$ws.Range("G5").Value = '9 - 12 Weeks'
$ws.Range("H5").Value = 'ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ'

# Row 6
$ws.Range("A6").Value = '''1327380'
$ws.Range("B6").Value = 'https://aiesec.org/opportunity/global-talent/1327380'
$ws.Range("C6").Value = 'Comunication Intern'
$ws.Range("D6").Value = 'Ümraniye, Elmalıkent, 34764 Ümraniye/İstanbul, Türkiye'
$ws.Range("E6").Value = 'No'
$ws.Range("F6").Value = '3 applicants'
$ws.Range("G6").Value = '9 - 12 Weeks'
$ws.Range("H6").Value = 'ENTES ELEKTRONİK CİHAZLAR İMALAT VE TİCARET ANONİM ŞİRKETİ'

# Row 7
$ws.Range("A7").Value = '''1327257'
$ws.Range("B7").Value = 'https://aiesec.org/opportunity/global-talent/1327257'
$ws.Range("C7").Value = 'Export Specialist at Arçek'
$ws.Range("D7").Value = 'Konya, Türkiye'
$ws.Range("E7").Value = 'No'
$ws.Range("F7").Value = '3 applicants'
$ws.Range("G7").Value = '6 - 18 Months'
$ws.Range("H7").Value = 'ARÇEK İNŞAAT TEKSTİL OTOMOTİV İÇ VE DIŞ TİCARET LİMİTED ŞİRK'

# Row 8
$ws.Range("A8").Value = '''1327143'
$ws.Range("B8").Value = 'https://aiesec.org/opportunity/global-talent/1327143'
$ws.Range("C8").Value = 'IT & Web Development Engineer'
$ws.Range("D8").Value = 'Ahmedabad, Gujarat, India'
$ws.Range("E8").Value = 'No'
$ws.Range("F8").Value = '0 applicants'
$ws.Range("G8").Value = '9 - 12 Weeks'
$ws.Range("H8").Value = 'WeHear Innovations Pvt Ltd'

# Row 9
$ws.Range("A9").Value = '''1327138'
$ws.Range("B9").Value = 'https://aiesec.org/opportunity/global-talent/1327138'
$ws.Range("C9").Value = 'UX Research Engineer'
$ws.Range("D9").Value = 'Ahmedabad, Gujarat, India'
$ws.Range("E9").Value = 'No'
$ws.Range("F9").Value = '0 applicants'
$ws.Range("G9").Value = '9 - 12 Weeks'
$ws.Range("H9").Value = 'WeHear Innovations Pvt Ltd'

# Row 10
$ws.Range("A10").Value = '''1325378'
$ws.Range("B10").Value = 'https://aiesec.org/opportunity/global-talent/1325378'
$ws.Range("C10").Value = 'Content Creation and Social Media Marketing Intern'
$ws.Range("D10").Value = 'Athens, Greece'
$ws.Range("E10").Value = 'No'
$ws.Range("F10").Value = '47 applicants'
$ws.Range("G10").Value = '9 - 12 Weeks'
$ws.Range("H10").Value = 'Eutopians'

# Row 11
$ws.Range("A11").Value = '''1314400'
$ws.Range("B11").Value = 'https://aiesec.org/opportunity/global-talent/1314400'
$ws.Range("C11").Value = '[Impact Curitiba] - Inside Sales'
$ws.Range("D11").Value = 'São Jorge D''Oeste - São Jorge d''Oeste, PR, 85575-000, Brasil'
$ws.Range("E11").Value = 'No'
$ws.Range("F11").Value = '45 applicants'
$ws.Range("G11").Value = '6 - 18 Months'
$ws.Range("H11").Value = 'Mocelin Indústria de Extintores'

# Remove now-obsolete rows 12-14 (table shrank from 14 to 11 data+header rows)
$ws.Range("A12:H14").EntireRow.Delete()

# Adjust column widths (C, D, F, H)
$ws.Columns.Item(3).ColumnWidth = 66.16666666666667
$ws.Columns.Item(4).ColumnWidth = 62.166666666666664
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 62.166666666666664

